$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts a new price observation as row 15 (same
# market/product/variety/quality as the existing row that is currently at
# row 15), pushing that row and everything below it down by one. Excel's
# row insert naturally carries the old row 42 down to row 43, which is
# exactly the new trailing row needed, so nothing else has to be appended.

$ws.Rows(15).Insert()

$newRow = @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44523, 10, "Fruta", 100103, "Frutos de hueso (carozo)", 100103003, "Damasco", "Castle Brite", "Segunda", 500, 28000, 28500, 28250, "`$/caja 18 kilos", "Provincia de Limarí", 1569, 18)

for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(15, $col).Value = $newRow[$col - 1]
}

# Column D carries the date number format; make sure the newly written
# value keeps it (Excel already propagates style on row insert, but set it
# explicitly in case the value write reset it).
$ws.Cells.Item(15, 4).NumberFormat = $ws.Cells.Item(16, 4).NumberFormat
